$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 12 (J12/K12) - empty cells formatted like the small divider row
# used elsewhere in the sheet. Build the format on J12 first (keeps the
# style table minimal), then copy that formatting across to K12.
# ---------------------------------------------------------------------------
$j12 = $ws.Range("J12")
$j12.Font.Size = 7
$j12.Font.Color = 0
$j12.Font.Name = "Courier New"
$j12.Font.Family = 3
$j12.HorizontalAlignment = -4131
$j12.VerticalAlignment = -4108

$j12.Copy()
$ws.Range("K12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# "RBC" column header label (shared string) for each of the per-group
# result tables.
# ---------------------------------------------------------------------------
$headerRows = @(4,7,10,15,18,21,24,29,32,35,38)
foreach ($r in $headerRows) {
    $ws.Range("M" + $r).Value = "RBC"
}

# ---------------------------------------------------------------------------
# p-values: doubled (two-tailed) test statistics.
# K8 / K11 keep a live formula; the rest are plain doubled values.
# ---------------------------------------------------------------------------
$ws.Range("K8").Formula = "=0.021558*2"
$ws.Range("K11").Formula = "=0.00124*2"

$ws.Range("K5").Value = 0.00019919532147312001
$ws.Range("K16").Value = 0.43707933762790901
$ws.Range("K19").Value = 0.002562
$ws.Range("K22").Value = 0.00021699999999999999
$ws.Range("K25").Value = 0.061078999208132402
$ws.Range("K30").Value = 0.00000223980479251966
$ws.Range("K33").Value = 0.0000000058949130000000004
$ws.Range("K36").Value = 0.00000000027044040000000002
$ws.Range("K39").Value = 0.000213983458902871

# ---------------------------------------------------------------------------
# New "RBC" effect-size values in column M.
# M5/M8/M11/M16/M22/M25/M33 are brand-new cells: copy the same bordered
# style used throughout column M/L ("style 1") from a donor cell before
# writing the value, so they end up styled just like their row neighbours.
# M19 intentionally stays unstyled (matches source). M30/M36/M39 already
# carry style 1 from the template, so only the value needs to be set.
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy()
$newStyledRows = @(5,8,11,16,22,25,33)
foreach ($r in $newStyledRows) {
    $ws.Range("M" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("M5").Value = 0.32100708103855202
$ws.Range("M8").Value = 0.174567
$ws.Range("M11").Value = 0.26111299999999998
$ws.Range("M16").Value = -0.064690026954177901
$ws.Range("M19").Value = 0.2807
$ws.Range("M22").Value = 0.344387
$ws.Range("M25").Value = 0.13834388516266499
$ws.Range("M30").Value = 0.375308641975308
$ws.Range("M33").Value = 0.46172800000000003
$ws.Range("M36").Value = 0.501017
$ws.Range("M39").Value = 0.27116920842411

# ---------------------------------------------------------------------------
# Sheet view: scroll down and move the selection.
# ---------------------------------------------------------------------------
try {
    $excel.ActiveWindow.ScrollRow = 12
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("O33").Select()
